$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert two new rows of fresh data at the top of the dataset
# (rows 789-790), shifting the existing rows 789-825 down to 791-827.
$ws.Rows("789:790").Insert()

# Row 789: new weekly price data (Sin especificar / Pintón)
$ws.Cells.Item(789, 1).Value = 4
$ws.Cells.Item(789, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(789, 3).Value = 'Los Lagos'
$ws.Cells.Item(789, 4).Value = 45041
$ws.Cells.Item(789, 5).Value = 10
$ws.Cells.Item(789, 6).Value = 'Fruta'
$ws.Cells.Item(789, 7).Value = 100108
$ws.Cells.Item(789, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(789, 9).Value = 100108006
$ws.Cells.Item(789, 10).Value = 'Plátano'
$ws.Cells.Item(789, 11).Value = 'Sin especificar'
$ws.Cells.Item(789, 12).Value = 'Pintón'
$ws.Cells.Item(789, 13).Value = 600
$ws.Cells.Item(789, 14).Value = 22000
$ws.Cells.Item(789, 15).Value = 22000
$ws.Cells.Item(789, 16).Value = 22000
$ws.Cells.Item(789, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(789, 18).Value = 'Ecuador'
$ws.Cells.Item(789, 19).Value = 1100
$ws.Cells.Item(789, 20).Value = 20
$ws.Cells.Item(789, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 790: new weekly price data (Sin especificar / Primera Pintón)
$ws.Cells.Item(790, 1).Value = 4
$ws.Cells.Item(790, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(790, 3).Value = 'Los Lagos'
$ws.Cells.Item(790, 4).Value = 45041
$ws.Cells.Item(790, 5).Value = 10
$ws.Cells.Item(790, 6).Value = 'Fruta'
$ws.Cells.Item(790, 7).Value = 100108
$ws.Cells.Item(790, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(790, 9).Value = 100108006
$ws.Cells.Item(790, 10).Value = 'Plátano'
$ws.Cells.Item(790, 11).Value = 'Sin especificar'
$ws.Cells.Item(790, 12).Value = 'Primera Pintón'
$ws.Cells.Item(790, 13).Value = 1200
$ws.Cells.Item(790, 14).Value = 23000
$ws.Cells.Item(790, 15).Value = 24000
$ws.Cells.Item(790, 16).Value = 23500
$ws.Cells.Item(790, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(790, 18).Value = 'Ecuador'
$ws.Cells.Item(790, 19).Value = 1175
$ws.Cells.Item(790, 20).Value = 20
$ws.Cells.Item(790, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 826: appended at the end of the dataset (duplicate of former row 824 content)
$ws.Cells.Item(826, 1).Value = 4
$ws.Cells.Item(826, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(826, 3).Value = 'Los Lagos'
$ws.Cells.Item(826, 4).Value = 44243
$ws.Cells.Item(826, 5).Value = 10
$ws.Cells.Item(826, 6).Value = 'Fruta'
$ws.Cells.Item(826, 7).Value = 100108
$ws.Cells.Item(826, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(826, 9).Value = 100108006
$ws.Cells.Item(826, 10).Value = 'Plátano'
$ws.Cells.Item(826, 11).Value = 'Sin especificar'
$ws.Cells.Item(826, 12).Value = 'Primera Pintón'
$ws.Cells.Item(826, 13).Value = 1100
$ws.Cells.Item(826, 14).Value = 17000
$ws.Cells.Item(826, 15).Value = 17500
$ws.Cells.Item(826, 16).Value = 17250
$ws.Cells.Item(826, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(826, 18).Value = 'Ecuador'
$ws.Cells.Item(826, 19).Value = 862
$ws.Cells.Item(826, 20).Value = 20
$ws.Cells.Item(826, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 827: appended at the end of the dataset (duplicate of former row 825 content)
$ws.Cells.Item(827, 1).Value = 4
$ws.Cells.Item(827, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(827, 3).Value = 'Los Lagos'
$ws.Cells.Item(827, 4).Value = 44217
$ws.Cells.Item(827, 5).Value = 10
$ws.Cells.Item(827, 6).Value = 'Fruta'
$ws.Cells.Item(827, 7).Value = 100108
$ws.Cells.Item(827, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(827, 9).Value = 100108006
$ws.Cells.Item(827, 10).Value = 'Plátano'
$ws.Cells.Item(827, 11).Value = 'Sin especificar'
$ws.Cells.Item(827, 12).Value = 'Primera Pintón'
$ws.Cells.Item(827, 13).Value = 800
$ws.Cells.Item(827, 14).Value = 17000
$ws.Cells.Item(827, 15).Value = 18000
$ws.Cells.Item(827, 16).Value = 17500
$ws.Cells.Item(827, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(827, 18).Value = 'Ecuador'
$ws.Cells.Item(827, 19).Value = 875
$ws.Cells.Item(827, 20).Value = 20
$ws.Cells.Item(827, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Refresh the used range / dimension so it reflects the two appended rows
$null = $ws.UsedRange
